# ---------------------------------------------------------------------------
# Adds a "Categoria" feature to the Estoque workbook:
#   * Insert a new "Categoria" column (E) on "Estoque" and "Removidos".
#   * Populate the category for every existing row.
#   * Add a brand-new "Categorias" sheet listing the available categories.
#   * Apply the rest of the row-level data changes (new items, edited items,
#     an item moved from Estoque to Removidos, ...).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Estoque sheet
# ---------------------------------------------------------------------------
$estoque = $wb.Worksheets.Item("Estoque")

# Insert the new "Categoria" column before the existing "Data_Criacao" column.
$estoque.Columns("E").Insert()
$estoque.Cells.Item(1, 5).Value = "Categoria"

# The item "cimento 20kg" (old row 8) was removed from stock (it now lives on
# the "Removidos" sheet) - delete it, shifting every row below up by one.
$estoque.Rows(8).Delete()

# Fill in the category for every item that stays on this sheet (rows 2-13
# after the delete above keep their original name/id/qty/price/dates).
$estoqueCategorias = @{
  2  = "cimento"
  3  = "maquinas"
  4  = "bucha"
  5  = "bucha"
  6  = "bucha"
  7  = "bucha"
  8  = "canos"
  9  = "canos"
  10 = "canos"
  11 = "canos"
  12 = "parafusos"
  13 = "prego"
}
foreach ($r in $estoqueCategorias.Keys) {
  $estoque.Cells.Item($r, 5).Value = $estoqueCategorias[$r]
}

# Refresh "Data_Alteracao" (column G) for the items above.
$estoqueAlteracao = @{
  2  = "29/10/2025 21:22"
  3  = "29/10/2025 21:24"
  4  = "29/10/2025 21:21"
  5  = "29/10/2025 21:21"
  6  = "29/10/2025 21:21"
  7  = "29/10/2025 21:21"
  8  = "29/10/2025 21:22"
  9  = "29/10/2025 21:22"
  10 = "29/10/2025 21:22"
  11 = "29/10/2025 21:22"
  12 = "29/10/2025 21:23"
  13 = "29/10/2025 21:23"
}
foreach ($r in $estoqueAlteracao.Keys) {
  $estoque.Cells.Item($r, 7).Value = $estoqueAlteracao[$r]
}

# Insert a brand-new item "cimento 20 kg" at row 14 (pushes "gugu" to row 15).
$estoque.Rows(14).Insert()
$estoque.Cells.Item(14, 1).Value = "cimento 20 kg"
$estoque.Cells.Item(14, 2).Value = "ID_15"
$estoque.Cells.Item(14, 3).Value = 15
$estoque.Cells.Item(14, 4).Value = 13
$estoque.Cells.Item(14, 5).Value = "cimento"
$estoque.Cells.Item(14, 6).Value = "29/10/2025 20:12"
$estoque.Cells.Item(14, 7).Value = "29/10/2025 20:12"

# Row 15 ("gugu"/ID_14) was edited into "maquita" (id and quantity unchanged).
$estoque.Cells.Item(15, 1).Value = "maquita"
$estoque.Cells.Item(15, 4).Value = 3000
$estoque.Cells.Item(15, 5).Value = "maquinas"
$estoque.Cells.Item(15, 6).Value = "29/10/2025 21:35"
$estoque.Cells.Item(15, 7).Value = "29/10/2025 21:35"

# Append a new "gugu" item as row 16.
$estoque.Cells.Item(16, 1).Value = "gugu"
$estoque.Cells.Item(16, 2).Value = "ID_15"
$estoque.Cells.Item(16, 3).Value = 2
$estoque.Cells.Item(16, 4).Value = 9999
$estoque.Cells.Item(16, 5).Value = "gugu"
$estoque.Cells.Item(16, 6).Value = "29/10/2025 21:37"
$estoque.Cells.Item(16, 7).Value = "29/10/2025 21:38"

# ---------------------------------------------------------------------------
# 2) Removidos sheet
# ---------------------------------------------------------------------------
$removidos = $wb.Worksheets.Item("Removidos")

# Insert the new "Categoria" column before the existing "Data_Criacao" column.
$removidos.Columns("E").Insert()
$removidos.Cells.Item(1, 5).Value = "Categoria"

# Every pre-existing removed item defaults to "Sem Categoria".
for ($r = 2; $r -le 18; $r++) {
  $removidos.Cells.Item($r, 5).Value = "Sem Categoria"
}

# New row 19: item "da" was created and removed the same day, with a category.
$removidos.Cells.Item(19, 1).Value = "da"
$removidos.Cells.Item(19, 2).Value = "ID_16"
$removidos.Cells.Item(19, 3).Value = 1
$removidos.Cells.Item(19, 5).Value = "cimento"
$removidos.Cells.Item(19, 6).Value = "29/10/2025 20:21"
$removidos.Cells.Item(19, 7).Value = "29/10/2025 20:21"
$removidos.Cells.Item(19, 8).Value = "29/10/2025 21:17"

# New row 20: "cimento 20kg" (previously Estoque row 8) got removed.
$removidos.Cells.Item(20, 1).Value = "cimento 20kg"
$removidos.Cells.Item(20, 2).Value = "ID_7"
$removidos.Cells.Item(20, 3).Value = 25
$removidos.Cells.Item(20, 4).Value = 406.9
$removidos.Cells.Item(20, 5).Value = "Sem Categoria"
$removidos.Cells.Item(20, 6).Value = "21/10/2025 22:29"
$removidos.Cells.Item(20, 7).Value = "21/10/2025 22:56"
$removidos.Cells.Item(20, 8).Value = "29/10/2025 21:21"

# New row 21: the original "gugu" (ID_14) item got removed.
$removidos.Cells.Item(21, 1).Value = "gugu"
$removidos.Cells.Item(21, 2).Value = "ID_14"
$removidos.Cells.Item(21, 3).Value = 1
$removidos.Cells.Item(21, 4).Value = 12
$removidos.Cells.Item(21, 5).Value = "Sem Categoria"
$removidos.Cells.Item(21, 6).Value = "23/10/2025 23:59"
$removidos.Cells.Item(21, 7).Value = "23/10/2025 23:59"
$removidos.Cells.Item(21, 8).Value = "29/10/2025 21:23"

# ---------------------------------------------------------------------------
# 3) New "Categorias" sheet
# ---------------------------------------------------------------------------
$categorias = $wb.Worksheets.Add($null, $removidos)
$categorias.Name = "Categorias"

# Copy the header formatting (bold, centered, bordered) from Estoque!A1.
$estoque.Range("A1").Copy($categorias.Range("A1"))
$categorias.Cells.Item(1, 1).Value = "Categoria"

$categoriaList = @(
  "Sem Categoria",
  "cimento",
  "bucha",
  "canos",
  "parafusos",
  "prego",
  "maquinas",
  "gugu"
)
$row = 2
foreach ($cat in $categoriaList) {
  $categorias.Cells.Item($row, 1).Value = $cat
  $row = $row + 1
}

# Restore the original active sheet/selection ("Estoque", cell A1).
$estoque.Activate() | Out-Null
$estoque.Range("A1").Select() | Out-Null

